# The author removed the two trailing "wrap-up" slides (HW reminder /
# open-lab instructions) from the deck - these correspond to
# sldId 276 (slide11.xml, creationId 2411685965) and
# sldId 277 (slide12.xml, creationId 2313341475), i.e. the last two
# slides in the deck (positions 11 and 12 of 12).
#
# Delete from the end first so the earlier deletion doesn't shift the
# index of the slide we still need to remove.
$p = $ppt.ActivePresentation

$last = $p.Slides.Count
$p.Slides.Item($last).Delete()
$p.Slides.Item($last - 1).Delete()
